$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update hours worked on the last logged day (45432) from 6.5 to 8.
# Downstream formulas (C64 running total, D2 = SUM(B:B), F2 = 40*D2)
# recalc automatically.
$ws.Range("B64").Value = 8

# Move the active selection from C64 to D64 to match the saved view state.
$ws.Range("D64").Select()
